$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow the cell writes below.
$ws.Unprotect()

# Update the confidential disclaimer text (A38) with the new "as of" date.
$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-14 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) figures for each holding row (2-35).
$ws.Range("D2").Value = 0.03571158590777979
$ws.Range("E2").Value = -0.0003944773175542649
$ws.Range("D3").Value = 0.0204100454060911
$ws.Range("E3").Value = -0.006264682850430692
$ws.Range("D4").Value = 0.01931051020089176
$ws.Range("E4").Value = -0.0008064516129030475
$ws.Range("D5").Value = 0.03727989891128229
$ws.Range("E5").Value = 0.004644515898535317
$ws.Range("D6").Value = 0.03445012232093515
$ws.Range("E6").Value = -0.000400000000000067
$ws.Range("D7").Value = 0.01983635481892061
$ws.Range("E7").Value = 0.00244625217896588
$ws.Range("D8").Value = 0.03700916449274528
$ws.Range("E8").Value = 0.01609096760351858
$ws.Range("D9").Value = 0.02040574234248521
$ws.Range("E9").Value = 0.005241753276095773
$ws.Range("D10").Value = 0.02638187805989829
$ws.Range("E10").Value = 0.0005825242718446866
$ws.Range("D11").Value = 0.02388456436005355
$ws.Range("E11").Value = 0.01474530831099186
$ws.Range("D12").Value = 0.05736383356838443
$ws.Range("E12").Value = 0.01072961373390569
$ws.Range("D13").Value = 0.02471187956523288
$ws.Range("E13").Value = 0.002611940298507642
$ws.Range("D14").Value = 0.02761060762716959
$ws.Range("E14").Value = -0.009432527004412217
$ws.Range("D15").Value = 0.03374175608828996
$ws.Range("E15").Value = -0.002040469307940884
$ws.Range("D16").Value = 0.01992574584263813
$ws.Range("E16").Value = -0.001995012468827828
$ws.Range("D17").Value = 0.03124843809036496
$ws.Range("E17").Value = 0.002213114754098466
$ws.Range("D18").Value = 0.04178428442158687
$ws.Range("E18").Value = 0.005590496156534064
$ws.Range("D19").Value = 0.1254149403253598
$ws.Range("E19").Value = 0.004694835680751019
$ws.Range("D20").Value = 0.009271411581411284
$ws.Range("E20").Value = -0.00508323802262034
$ws.Range("D21").Value = 0.01542863455890541
$ws.Range("E21").Value = 0.004601868637568174
$ws.Range("D22").Value = 0.0171123618755515
$ws.Range("E22").Value = 0.007430025445292587
$ws.Range("D23").Value = 0.01541111494279573
$ws.Range("E23").Value = 0.004710144927536186
$ws.Range("D24").Value = 0.02120375573358637
$ws.Range("E24").Value = 0.01111572823602724
$ws.Range("D25").Value = 0.01252411785187983
$ws.Range("E25").Value = 0.01562071800493281
$ws.Range("D26").Value = 0.0424919847039121
$ws.Range("E26").Value = 0.00569630553897893
$ws.Range("D27").Value = 0.02402390165776796
$ws.Range("E27").Value = 0.0001961745953900085
$ws.Range("D28").Value = 0.04541299648809358
$ws.Range("E28").Value = 0.005285920230658325
$ws.Range("D29").Value = 0.05522345284741424
$ws.Range("E29").Value = 0.01671248418204829
$ws.Range("D30").Value = 0.01280653201448997
$ws.Range("E30").Value = 0.03253652058432932
$ws.Range("D31").Value = 0.02064445991871652
$ws.Range("E31").Value = 0.002692307692307772
$ws.Range("D32").Value = 0.01324856934609941
$ws.Range("E32").Value = 0.006283228612856462
$ws.Range("D33").Value = 0.04197285081603053
$ws.Range("E33").Value = 0.003613835828601086
$ws.Range("D34").Value = 0.01674250331323604
$ws.Range("E34").Value = 0.01844077961019486
$ws.Range("E35").Value = 0.005403550095504528

# Restore sheet protection.
$ws.Protect()
